$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build strings containing special/unicode characters explicitly (avoids
# PowerShell escape-sequence support gaps in this runtime).
$ldq = [char]0x201C   # “
$rdq = [char]0x201D   # ”

# Row 6 (EXE-02 scenario "Register using already registered email") --
# update the test-description / results / status / severity / bug-id
# columns to reflect the new (failed) test run.
$ws.Range("I6").Value = "Register using already registered email"
$ws.Range("J6").Value = "Error message appears"
$ws.Range("K6").Value = "Error message " + $ldq + "Email already registered" + $rdq + " displayed"
$ws.Range("L6").Value = "FAIL"
$ws.Range("M6").Value = "LOW"
$ws.Range("N6").Value = "BUG-01"

# Test date moved back a day.
$ws.Range("P6").Value = 46076

# L6 gets a new highlight style: same font/border/alignment as before,
# but filled yellow instead of the green "PASS" fill.
$ws.Range("L6").Interior.Color = 65535

# Column K now holds much longer text, so widen it independently of the
# rest of the K:O block (which used to share one <col> run).
$ws.Columns("K").ColumnWidth = 45.5

# Update the view: scroll so column D is leftmost and select N6.
$ws.Range("N6").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
